# Apply rerun simulation data + two new HKL entries (Holden, Rizzie Spiral)
# and rename "Thomas Hex" -> "Matthies Hex".
# Rows 4-29 (A=2..27) take the values that used to live two rows higher
# (old rows 4-27), rows 4 and 5 get freshly computed values for the two new
# HKL entries, and two new rows (30, 31 / A=28,29) are appended for the last
# two original entries that got pushed off the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = "Holden"
$ws.Cells.Item(4, 3).Value = 1.011964140332335
$ws.Cells.Item(4, 4).Value = 0.9902244328948805
$ws.Cells.Item(4, 5).Value = 1.00434370967539
$ws.Cells.Item(4, 6).Value = 1.004438160217789
$ws.Cells.Item(4, 7).Value = 0.9822823437521215
$ws.Cells.Item(4, 8).Value = 1.004477025573838
$ws.Cells.Item(4, 9).Value = 1.00434370967539
$ws.Cells.Item(4, 10).Value = 0.9822823437521215
$ws.Cells.Item(4, 11).Value = 1.00434370967539
$ws.Cells.Item(4, 12).Value = 1.004477025573838
$ws.Cells.Item(4, 13).Value = 0.9933796846629798
$ws.Cells.Item(4, 14).Value = 0.9933796846629798
$ws.Cells.Item(4, 15).Value = 0.9923279340736134
$ws.Cells.Item(4, 16).Value = 0.9970343596671167
$ws.Cells.Item(4, 17).Value = 0.9970343596671167
$ws.Cells.Item(4, 18).Value = 0.9988616971691852
$ws.Cells.Item(4, 19).Value = 0.9988616971691852
$ws.Cells.Item(4, 20).Value = 0.9996216354077258
# Row 5
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$ws.Cells.Item(5, 3).Value = 1.008455316609172
$ws.Cells.Item(5, 4).Value = 0.9932676715771742
$ws.Cells.Item(5, 5).Value = 1.00280431016725
$ws.Cells.Item(5, 6).Value = 1.003046604259772
$ws.Cells.Item(5, 7).Value = 0.9879033132086472
$ws.Cells.Item(5, 8).Value = 1.003146314203467
$ws.Cells.Item(5, 9).Value = 1.00280431016725
$ws.Cells.Item(5, 10).Value = 0.9879033132086472
$ws.Cells.Item(5, 11).Value = 1.00280431016725
$ws.Cells.Item(5, 12).Value = 1.003146314203467
$ws.Cells.Item(5, 13).Value = 0.9955248137060572
$ws.Cells.Item(5, 14).Value = 0.9955248137060572
$ws.Cells.Item(5, 15).Value = 0.9947724329964295
$ws.Cells.Item(5, 16).Value = 0.9979513125264549
$ws.Cells.Item(5, 17).Value = 0.9979513125264549
$ws.Cells.Item(5, 18).Value = 0.9991645619366536
$ws.Cells.Item(5, 19).Value = 0.9991645619366536
$ws.Cells.Item(5, 20).Value = 0.9997705883375803
# Row 6
$ws.Cells.Item(6, 2).Value = "RotRing OmegaMax-90"
$ws.Cells.Item(6, 3).Value = 1.00434741482876
$ws.Cells.Item(6, 4).Value = 0.9963343513138099
$ws.Cells.Item(6, 5).Value = 1.001749302883786
$ws.Cells.Item(6, 6).Value = 1.001670605027955
$ws.Cells.Item(6, 7).Value = 0.9932884434558865
$ws.Cells.Item(6, 8).Value = 1.001638212689208
$ws.Cells.Item(6, 9).Value = 1.001749302883786
$ws.Cells.Item(6, 10).Value = 0.9932884434558865
$ws.Cells.Item(6, 11).Value = 1.001749302883786
$ws.Cells.Item(6, 12).Value = 1.001638212689208
$ws.Cells.Item(6, 13).Value = 0.9974633280725471
$ws.Cells.Item(6, 14).Value = 0.9974633280725471
$ws.Cells.Item(6, 15).Value = 0.9970870024863014
$ws.Cells.Item(6, 16).Value = 0.9988919863429601
$ws.Cells.Item(6, 17).Value = 0.9988919863429601
$ws.Cells.Item(6, 18).Value = 0.9996063154781665
$ws.Cells.Item(6, 19).Value = 0.9996063154781665
$ws.Cells.Item(6, 20).Value = 0.9998380550332343
# Row 7
$ws.Cells.Item(7, 2).Value = "Equal Angle"
$ws.Cells.Item(7, 3).Value = 1.014260901224783
$ws.Cells.Item(7, 4).Value = 0.9883449372262256
$ws.Cells.Item(7, 5).Value = 1.005181905914989
$ws.Cells.Item(7, 6).Value = 1.005291622146976
$ws.Cells.Item(7, 7).Value = 0.978874146433716
$ws.Cells.Item(7, 8).Value = 1.005336770626802
$ws.Cells.Item(7, 9).Value = 1.005181905914989
$ws.Cells.Item(7, 10).Value = 0.978874146433716
$ws.Cells.Item(7, 11).Value = 1.005181905914989
$ws.Cells.Item(7, 12).Value = 1.005336770626802
$ws.Cells.Item(7, 13).Value = 0.9921054585302589
$ws.Cells.Item(7, 14).Value = 0.9921054585302589
$ws.Cells.Item(7, 15).Value = 0.9908519514289145
$ws.Cells.Item(7, 16).Value = 0.9964642743251689
$ws.Cells.Item(7, 17).Value = 0.9964642743251689
$ws.Cells.Item(7, 18).Value = 0.9986436822226239
$ws.Cells.Item(7, 19).Value = 0.9986436822226239
$ws.Cells.Item(7, 20).Value = 0.9995483805955819
# Row 8
$ws.Cells.Item(8, 2).Value = "Tilt Rotate"
$ws.Cells.Item(8, 3).Value = 1.048028837802456
$ws.Cells.Item(8, 4).Value = 0.9609071189023591
$ws.Cells.Item(8, 5).Value = 1.017211286608418
$ws.Cells.Item(8, 6).Value = 1.017739930091168
$ws.Cells.Item(8, 7).Value = 0.9292360345761207
$ws.Cells.Item(8, 8).Value = 1.017957487995727
$ws.Cells.Item(8, 9).Value = 1.017211286608418
$ws.Cells.Item(8, 10).Value = 0.9292360345761207
$ws.Cells.Item(8, 11).Value = 1.017211286608418
$ws.Cells.Item(8, 12).Value = 1.017957487995727
$ws.Cells.Item(8, 13).Value = 0.973596761285924
$ws.Cells.Item(8, 14).Value = 0.973596761285924
$ws.Cells.Item(8, 15).Value = 0.9693668804914024
$ws.Cells.Item(8, 16).Value = 0.9881349363934219
$ws.Cells.Item(8, 17).Value = 0.9881349363934219
$ws.Cells.Item(8, 18).Value = 0.9954040239471708
$ws.Cells.Item(8, 19).Value = 0.9954040239471708
$ws.Cells.Item(8, 20).Value = 0.9985134493293749
# Row 9
$ws.Cells.Item(9, 2).Value = "CLR"
$ws.Cells.Item(9, 3).Value = 1.001460883193413
$ws.Cells.Item(9, 4).Value = 0.9987829459745278
$ws.Cells.Item(9, 5).Value = 1.000565638231457
$ws.Cells.Item(9, 6).Value = 1.000553864583935
$ws.Cells.Item(9, 7).Value = 0.9977801839470487
$ws.Cells.Item(9, 8).Value = 1.000549015190798
$ws.Cells.Item(9, 9).Value = 1.000565638231457
$ws.Cells.Item(9, 10).Value = 0.9977801839470487
$ws.Cells.Item(9, 11).Value = 1.000565638231457
$ws.Cells.Item(9, 12).Value = 1.000549015190798
$ws.Cells.Item(9, 13).Value = 0.9991645995689236
$ws.Cells.Item(9, 14).Value = 0.9991645995689236
$ws.Cells.Item(9, 15).Value = 0.999037381704125
$ws.Cells.Item(9, 16).Value = 0.9996316124564345
$ws.Cells.Item(9, 17).Value = 0.9996316124564345
$ws.Cells.Item(9, 18).Value = 0.99986511890019
$ws.Cells.Item(9, 19).Value = 0.99986511890019
$ws.Cells.Item(9, 20).Value = 0.9999487551868634
# Row 10
$ws.Cells.Item(10, 2).Value = "Rizzie Hex"
$ws.Cells.Item(10, 3).Value = 1.000122421706267
$ws.Cells.Item(10, 4).Value = 0.9999026844663654
$ws.Cells.Item(10, 5).Value = 1.00004036131163
$ws.Cells.Item(10, 6).Value = 1.000044029391928
$ws.Cells.Item(10, 7).Value = 0.9998252440252546
$ws.Cells.Item(10, 8).Value = 1.000045535591986
$ws.Cells.Item(10, 9).Value = 1.00004036131163
$ws.Cells.Item(10, 10).Value = 0.9998252440252546
$ws.Cells.Item(10, 11).Value = 1.00004036131163
$ws.Cells.Item(10, 12).Value = 1.000045535591986
$ws.Cells.Item(10, 13).Value = 0.9999353898086204
$ws.Cells.Item(10, 14).Value = 0.9999353898086204
$ws.Cells.Item(10, 15).Value = 0.9999244880278688
$ws.Cells.Item(10, 16).Value = 0.9999703803096237
$ws.Cells.Item(10, 17).Value = 0.9999703803096237
$ws.Cells.Item(10, 18).Value = 0.9999878755601254
$ws.Cells.Item(10, 19).Value = 0.9999878755601254
$ws.Cells.Item(10, 20).Value = 0.9999967127489052
# Row 11
$ws.Cells.Item(11, 2).Value = "Matthies Hex"
$ws.Cells.Item(11, 3).Value = 1.002581527818124
$ws.Cells.Item(11, 4).Value = 0.9978480652736793
$ws.Cells.Item(11, 5).Value = 1.001001467538597
$ws.Cells.Item(11, 6).Value = 1.00097938565425
$ws.Cells.Item(11, 7).Value = 0.9960742767011835
$ws.Cells.Item(11, 8).Value = 1.000970296424139
$ws.Cells.Item(11, 9).Value = 1.001001467538597
$ws.Cells.Item(11, 10).Value = 0.9960742767011835
$ws.Cells.Item(11, 11).Value = 1.001001467538597
$ws.Cells.Item(11, 12).Value = 1.000970296424139
$ws.Cells.Item(11, 13).Value = 0.9985222865626615
$ws.Cells.Item(11, 14).Value = 0.9985222865626615
$ws.Cells.Item(11, 15).Value = 0.9982975461330007
$ws.Cells.Item(11, 16).Value = 0.9993486802213066
$ws.Cells.Item(11, 17).Value = 0.9993486802213066
$ws.Cells.Item(11, 18).Value = 0.9997618770506291
$ws.Cells.Item(11, 19).Value = 0.9997618770506291
$ws.Cells.Item(11, 20).Value = 0.9999091699016621
# Row 12
$ws.Cells.Item(12, 2).Value = "Tilt Rotate_Partial"
$ws.Cells.Item(12, 3).Value = 1.048753775722664
$ws.Cells.Item(12, 4).Value = 0.9603255909897953
$ws.Cells.Item(12, 5).Value = 1.017458220011402
$ws.Cells.Item(12, 6).Value = 1.018003339677607
$ws.Cells.Item(12, 7).Value = 0.9281884986228039
$ws.Cells.Item(12, 8).Value = 1.018227678764699
$ws.Cells.Item(12, 9).Value = 1.017458220011402
$ws.Cells.Item(12, 10).Value = 0.9281884986228039
$ws.Cells.Item(12, 11).Value = 1.017458220011402
$ws.Cells.Item(12, 12).Value = 1.018227678764699
$ws.Cells.Item(12, 13).Value = 0.9732080886937513
$ws.Cells.Item(12, 14).Value = 0.9732080886937513
$ws.Cells.Item(12, 15).Value = 0.9689139227924327
$ws.Cells.Item(12, 16).Value = 0.9879581324663015
$ws.Cells.Item(12, 17).Value = 0.9879581324663015
$ws.Cells.Item(12, 18).Value = 0.9953331543525767
$ws.Cells.Item(12, 19).Value = 0.9953331543525767
$ws.Cells.Item(12, 20).Value = 0.9984928506314952
# Row 13
$ws.Cells.Item(13, 2).Value = "RotRing OmegaMax-60"
$ws.Cells.Item(13, 3).Value = 1.004121746080861
$ws.Cells.Item(13, 4).Value = 0.9965451446109521
$ws.Cells.Item(13, 5).Value = 1.001627604005186
$ws.Cells.Item(13, 6).Value = 1.001573419343117
$ws.Cells.Item(13, 7).Value = 0.9936862626508077
$ws.Cells.Item(13, 8).Value = 1.001551114738082
$ws.Cells.Item(13, 9).Value = 1.001627604005186
$ws.Cells.Item(13, 10).Value = 0.9936862626508077
$ws.Cells.Item(13, 11).Value = 1.001627604005186
$ws.Cells.Item(13, 12).Value = 1.001551114738082
$ws.Cells.Item(13, 13).Value = 0.9976186886944447
$ws.Cells.Item(13, 14).Value = 0.9976186886944447
$ws.Cells.Item(13, 15).Value = 0.9972608406666138
$ws.Cells.Item(13, 16).Value = 0.9989549937980251
$ws.Cells.Item(13, 17).Value = 0.998954993798025
$ws.Cells.Item(13, 18).Value = 0.9996231463498152
$ws.Cells.Item(13, 19).Value = 0.9996231463498152
$ws.Cells.Item(13, 20).Value = 0.999850881904834
# Row 14
$ws.Cells.Item(14, 2).Value = "Equal Angle_Partial"
$ws.Cells.Item(14, 3).Value = 1.014635376810525
$ws.Cells.Item(14, 4).Value = 0.9880943623578948
$ws.Cells.Item(14, 5).Value = 1.005234431284211
$ws.Cells.Item(14, 6).Value = 1.005402270168422
$ws.Cells.Item(14, 7).Value = 0.9784530808526318
$ws.Cells.Item(14, 8).Value = 1.005471339473684
$ws.Cells.Item(14, 9).Value = 1.005234431284211
$ws.Cells.Item(14, 10).Value = 0.9784530808526318
$ws.Cells.Item(14, 11).Value = 1.005234431284211
$ws.Cells.Item(14, 12).Value = 1.005471339473684
$ws.Cells.Item(14, 13).Value = 0.9919622101631578
$ws.Cells.Item(14, 14).Value = 0.9919622101631578
$ws.Cells.Item(14, 15).Value = 0.9906729275614033
$ws.Cells.Item(14, 16).Value = 0.9963862838701756
$ws.Cells.Item(14, 17).Value = 0.9963862838701756
$ws.Cells.Item(14, 18).Value = 0.9985983207236845
$ws.Cells.Item(14, 19).Value = 0.9985983207236845
$ws.Cells.Item(14, 20).Value = 0.9995484768245615
# Row 15
$ws.Cells.Item(15, 2).Value = "Rizzie Hex_Partial"
$ws.Cells.Item(15, 3).Value = 0.9905104659202076
$ws.Cells.Item(15, 4).Value = 1.007805984251172
$ws.Cells.Item(15, 5).Value = 0.9964758800235556
$ws.Cells.Item(15, 6).Value = 0.9964530973710748
$ws.Cells.Item(15, 7).Value = 1.014179175022297
$ws.Cells.Item(15, 8).Value = 0.9964437224966527
$ws.Cells.Item(15, 9).Value = 0.9964758800235556
$ws.Cells.Item(15, 10).Value = 1.014179175022297
$ws.Cells.Item(15, 11).Value = 0.9964758800235556
$ws.Cells.Item(15, 12).Value = 0.9964437224966527
$ws.Cells.Item(15, 13).Value = 1.005311448759475
$ws.Cells.Item(15, 14).Value = 1.005311448759475
$ws.Cells.Item(15, 15).Value = 1.00614296059004
$ws.Cells.Item(15, 16).Value = 1.002366259180835
$ws.Cells.Item(15, 17).Value = 1.002366259180835
$ws.Cells.Item(15, 18).Value = 1.000893664391515
$ws.Cells.Item(15, 19).Value = 1.000893664391515
$ws.Cells.Item(15, 20).Value = 1.00031138751416
# Row 16
$ws.Cells.Item(16, 2).Value = "ND Single"
$ws.Cells.Item(16, 3).Value = 1.084257500000001
$ws.Cells.Item(16, 4).Value = 0.9314548999999992
$ws.Cells.Item(16, 5).Value = 1.0301398
$ws.Cells.Item(16, 6).Value = 1.031103000000001
$ws.Cells.Item(16, 7).Value = 0.8759446600000014
$ws.Cells.Item(16, 8).Value = 1.031499400000001
$ws.Cells.Item(16, 9).Value = 1.0301398
$ws.Cells.Item(16, 10).Value = 0.8759446600000014
$ws.Cells.Item(16, 11).Value = 1.0301398
$ws.Cells.Item(16, 12).Value = 1.031499400000001
$ws.Cells.Item(16, 13).Value = 0.9537220300000011
$ws.Cells.Item(16, 14).Value = 0.9537220300000011
$ws.Cells.Item(16, 15).Value = 0.9462996533333339
$ws.Cells.Item(16, 16).Value = 0.9791946200000009
$ws.Cells.Item(16, 17).Value = 0.9791946200000009
$ws.Cells.Item(16, 18).Value = 0.9919309150000009
$ws.Cells.Item(16, 19).Value = 0.9919309150000009
$ws.Cells.Item(16, 20).Value = 0.9973998766666675
# Row 17
$ws.Cells.Item(17, 2).Value = "RD Single"
$ws.Cells.Item(17, 3).Value = 1.0315004
$ws.Cells.Item(17, 4).Value = 0.97429176
$ws.Cells.Item(17, 5).Value = 1.0113915
$ws.Cells.Item(17, 6).Value = 1.01167
$ws.Cells.Item(17, 7).Value = 0.95342312
$ws.Cells.Item(17, 8).Value = 1.0117846
$ws.Cells.Item(17, 9).Value = 1.0113915
$ws.Cells.Item(17, 10).Value = 0.95342312
$ws.Cells.Item(17, 11).Value = 1.0113915
$ws.Cells.Item(17, 12).Value = 1.0117846
$ws.Cells.Item(17, 13).Value = 0.98260386
$ws.Cells.Item(17, 14).Value = 0.98260386
$ws.Cells.Item(17, 15).Value = 0.97983316
$ws.Cells.Item(17, 16).Value = 0.9921997399999999
$ws.Cells.Item(17, 17).Value = 0.9921997399999999
$ws.Cells.Item(17, 18).Value = 0.99699768
$ws.Cells.Item(17, 19).Value = 0.99699768
$ws.Cells.Item(17, 20).Value = 0.99901023
# Row 18
$ws.Cells.Item(18, 2).Value = "TD Single"
$ws.Cells.Item(18, 3).Value = 1.0301398
$ws.Cells.Item(18, 4).Value = 0.9742505199999999
$ws.Cells.Item(18, 5).Value = 1.012634
$ws.Cells.Item(18, 6).Value = 1.0117535
$ws.Cells.Item(18, 7).Value = 0.9526597999999999
$ws.Cells.Item(18, 8).Value = 1.0113912
$ws.Cells.Item(18, 9).Value = 1.012634
$ws.Cells.Item(18, 10).Value = 0.9526597999999999
$ws.Cells.Item(18, 11).Value = 1.012634
$ws.Cells.Item(18, 12).Value = 1.0113912
$ws.Cells.Item(18, 13).Value = 0.9820255
$ws.Cells.Item(18, 14).Value = 0.9820255
$ws.Cells.Item(18, 15).Value = 0.97943384
$ws.Cells.Item(18, 16).Value = 0.9922283333333333
$ws.Cells.Item(18, 17).Value = 0.9922283333333333
$ws.Cells.Item(18, 18).Value = 0.99732975
$ws.Cells.Item(18, 19).Value = 0.99732975
$ws.Cells.Item(18, 20).Value = 0.9988048033333333
# Row 19
$ws.Cells.Item(19, 2).Value = "Morris Single"
$ws.Cells.Item(19, 3).Value = 0.9760500600000001
$ws.Cells.Item(19, 4).Value = 1.0198259
$ws.Cells.Item(19, 5).Value = 0.99091759
$ws.Cells.Item(19, 6).Value = 0.99098449
$ws.Cells.Item(19, 7).Value = 1.0360868
$ws.Cells.Item(19, 8).Value = 0.9910120299999999
$ws.Cells.Item(19, 9).Value = 0.99091759
$ws.Cells.Item(19, 10).Value = 1.0360868
$ws.Cells.Item(19, 11).Value = 0.99091759
$ws.Cells.Item(19, 12).Value = 0.9910120299999999
$ws.Cells.Item(19, 13).Value = 1.013549415
$ws.Cells.Item(19, 14).Value = 1.013549415
$ws.Cells.Item(19, 15).Value = 1.015641576666667
$ws.Cells.Item(19, 16).Value = 1.006005473333333
$ws.Cells.Item(19, 17).Value = 1.006005473333333
$ws.Cells.Item(19, 18).Value = 1.0022335025
$ws.Cells.Item(19, 19).Value = 1.0022335025
$ws.Cells.Item(19, 20).Value = 1.000812811666667
# Row 20
$ws.Cells.Item(20, 2).Value = "Ring Perpendicular to ND"
$ws.Cells.Item(20, 3).Value = 1.030828920547945
$ws.Cells.Item(20, 4).Value = 0.9742718257534243
$ws.Cells.Item(20, 5).Value = 1.012004083561644
$ws.Cells.Item(20, 6).Value = 1.011710994520548
$ws.Cells.Item(20, 7).Value = 0.9530474284931509
$ws.Cells.Item(20, 8).Value = 1.01159038630137
$ws.Cells.Item(20, 9).Value = 1.012004083561644
$ws.Cells.Item(20, 10).Value = 0.9530474284931509
$ws.Cells.Item(20, 11).Value = 1.012004083561644
$ws.Cells.Item(20, 12).Value = 1.01159038630137
$ws.Cells.Item(20, 13).Value = 0.9823189073972604
$ws.Cells.Item(20, 14).Value = 0.9823189073972604
$ws.Cells.Item(20, 15).Value = 0.979636546849315
$ws.Cells.Item(20, 16).Value = 0.9922139661187216
$ws.Cells.Item(20, 17).Value = 0.9922139661187216
$ws.Cells.Item(20, 18).Value = 0.9971614954794523
$ws.Cells.Item(20, 19).Value = 0.9971614954794523
$ws.Cells.Item(20, 20).Value = 0.9989089398630138
# Row 21
$ws.Cells.Item(21, 2).Value = "Ring Perpendicular to RD"
$ws.Cells.Item(21, 3).Value = 1.013988627368421
$ws.Cells.Item(21, 4).Value = 0.9883787936842104
$ws.Cells.Item(21, 5).Value = 1.005367114736842
$ws.Cells.Item(21, 6).Value = 1.005286852105263
$ws.Cells.Item(21, 7).Value = 0.9788228542105264
$ws.Cells.Item(21, 8).Value = 1.005253827894737
$ws.Cells.Item(21, 9).Value = 1.005367114736842
$ws.Cells.Item(21, 10).Value = 0.9788228542105264
$ws.Cells.Item(21, 11).Value = 1.005367114736842
$ws.Cells.Item(21, 12).Value = 1.005253827894737
$ws.Cells.Item(21, 13).Value = 0.9920383410526317
$ws.Cells.Item(21, 14).Value = 0.9920383410526317
$ws.Cells.Item(21, 15).Value = 0.9908184919298245
$ws.Cells.Item(21, 16).Value = 0.9964812656140353
$ws.Cells.Item(21, 17).Value = 0.9964812656140353
$ws.Cells.Item(21, 18).Value = 0.9987027278947369
$ws.Cells.Item(21, 19).Value = 0.9987027278947369
$ws.Cells.Item(21, 20).Value = 0.999516345
# Row 22
$ws.Cells.Item(22, 2).Value = "Ring Perpendicular to TD"
$ws.Cells.Item(22, 3).Value = 1.014426614736842
$ws.Cells.Item(22, 4).Value = 0.9883926205263159
$ws.Cells.Item(22, 5).Value = 1.004966334210526
$ws.Cells.Item(22, 6).Value = 1.005259682105263
$ws.Cells.Item(22, 7).Value = 0.9790699068421053
$ws.Cells.Item(22, 8).Value = 1.005380401578948
$ws.Cells.Item(22, 9).Value = 1.004966334210526
$ws.Cells.Item(22, 10).Value = 0.9790699068421053
$ws.Cells.Item(22, 11).Value = 1.004966334210526
$ws.Cells.Item(22, 12).Value = 1.005380401578948
$ws.Cells.Item(22, 13).Value = 0.9922251542105265
$ws.Cells.Item(22, 14).Value = 0.9922251542105265
$ws.Cells.Item(22, 15).Value = 0.9909476429824563
$ws.Cells.Item(22, 16).Value = 0.9964722142105265
$ws.Cells.Item(22, 17).Value = 0.9964722142105265
$ws.Cells.Item(22, 18).Value = 0.9985957442105264
$ws.Cells.Item(22, 19).Value = 0.9985957442105264
$ws.Cells.Item(22, 20).Value = 0.9995825933333334
# Row 23
$ws.Cells.Item(23, 2).Value = "OffsetFTD"
$ws.Cells.Item(23, 3).Value = 0.9888323744684009
$ws.Cells.Item(23, 4).Value = 1.00897834241764
$ws.Cells.Item(23, 5).Value = 0.9961659865889145
$ws.Cells.Item(23, 6).Value = 0.9959320169001692
$ws.Cells.Item(23, 7).Value = 1.016185282777185
$ws.Cells.Item(23, 8).Value = 0.9958357270684536
$ws.Cells.Item(23, 9).Value = 0.9961659865889145
$ws.Cells.Item(23, 10).Value = 1.016185282777185
$ws.Cells.Item(23, 11).Value = 0.9961659865889145
$ws.Cells.Item(23, 12).Value = 0.9958357270684536
$ws.Cells.Item(23, 13).Value = 1.006010504922819
$ws.Cells.Item(23, 14).Value = 1.006010504922819
$ws.Cells.Item(23, 15).Value = 1.006999784087759
$ws.Cells.Item(23, 16).Value = 1.002728998811518
$ws.Cells.Item(23, 17).Value = 1.002728998811518
$ws.Cells.Item(23, 18).Value = 1.001088245755867
$ws.Cells.Item(23, 19).Value = 1.001088245755867
$ws.Cells.Item(23, 20).Value = 1.000321621703461
# Row 24
$ws.Cells.Item(24, 2).Value = "OffsetATD"
$ws.Cells.Item(24, 3).Value = 0.99730459539724
$ws.Cells.Item(24, 4).Value = 1.002323897636629
$ws.Cells.Item(24, 5).Value = 0.9988383306202269
$ws.Cells.Item(24, 6).Value = 0.9989381065659431
$ws.Cells.Item(24, 7).Value = 1.004284544148603
$ws.Cells.Item(24, 8).Value = 0.9989791665705556
$ws.Cells.Item(24, 9).Value = 0.9988383306202269
$ws.Cells.Item(24, 10).Value = 1.004284544148603
$ws.Cells.Item(24, 11).Value = 0.9988383306202269
$ws.Cells.Item(24, 12).Value = 0.9989791665705556
$ws.Cells.Item(24, 13).Value = 1.001631855359579
$ws.Cells.Item(24, 14).Value = 1.001631855359579
$ws.Cells.Item(24, 15).Value = 1.001862536118596
$ws.Cells.Item(24, 16).Value = 1.000700680446462
$ws.Cells.Item(24, 17).Value = 1.000700680446462
$ws.Cells.Item(24, 18).Value = 1.000235092989903
$ws.Cells.Item(24, 19).Value = 1.000235092989903
$ws.Cells.Item(24, 20).Value = 1.000111440156533
# Row 25
$ws.Cells.Item(25, 2).Value = "OffsetF45"
$ws.Cells.Item(25, 3).Value = 0.9888456298233126
$ws.Cells.Item(25, 4).Value = 1.008979219520659
$ws.Cells.Item(25, 5).Value = 0.996153154398043
$ws.Cells.Item(25, 6).Value = 0.9959309568614352
$ws.Cells.Item(25, 7).Value = 1.016193888586086
$ws.Cells.Item(25, 8).Value = 0.995839515633447
$ws.Cells.Item(25, 9).Value = 0.996153154398043
$ws.Cells.Item(25, 10).Value = 1.016193888586086
$ws.Cells.Item(25, 11).Value = 0.996153154398043
$ws.Cells.Item(25, 12).Value = 0.995839515633447
$ws.Cells.Item(25, 13).Value = 1.006016702109766
$ws.Cells.Item(25, 14).Value = 1.006016702109766
$ws.Cells.Item(25, 15).Value = 1.007004207913397
$ws.Cells.Item(25, 16).Value = 1.002728852872525
$ws.Cells.Item(25, 17).Value = 1.002728852872525
$ws.Cells.Item(25, 18).Value = 1.001084928253905
$ws.Cells.Item(25, 19).Value = 1.001084928253905
$ws.Cells.Item(25, 20).Value = 1.000323727470497
# Row 26
$ws.Cells.Item(26, 2).Value = "OffsetA45"
$ws.Cells.Item(26, 3).Value = 0.9972987008390212
$ws.Cells.Item(26, 4).Value = 1.002323846036044
$ws.Cells.Item(26, 5).Value = 0.9988435240341823
$ws.Cells.Item(26, 6).Value = 0.998938403644375
$ws.Cells.Item(26, 7).Value = 1.004281537449279
$ws.Cells.Item(26, 8).Value = 0.9989774478140858
$ws.Cells.Item(26, 9).Value = 0.9988435240341823
$ws.Cells.Item(26, 10).Value = 1.004281537449279
$ws.Cells.Item(26, 11).Value = 0.9988435240341823
$ws.Cells.Item(26, 12).Value = 0.9989774478140858
$ws.Cells.Item(26, 13).Value = 1.001629492631682
$ws.Cells.Item(26, 14).Value = 1.001629492631682
$ws.Cells.Item(26, 15).Value = 1.001860943766469
$ws.Cells.Item(26, 16).Value = 1.000700836432516
$ws.Cells.Item(26, 17).Value = 1.000700836432516
$ws.Cells.Item(26, 18).Value = 1.000236508332932
$ws.Cells.Item(26, 19).Value = 1.000236508332932
$ws.Cells.Item(26, 20).Value = 1.000110576636164
# Row 27
$ws.Cells.Item(27, 2).Value = "OffsetFRD"
$ws.Cells.Item(27, 3).Value = 0.9888590492768324
$ws.Cells.Item(27, 4).Value = 1.00897997056303
$ws.Cells.Item(27, 5).Value = 0.9961403935434965
$ws.Cells.Item(27, 6).Value = 0.9959299642086804
$ws.Cells.Item(27, 7).Value = 1.016202236447099
$ws.Cells.Item(27, 8).Value = 0.9958433541799609
$ws.Cells.Item(27, 9).Value = 0.9961403935434965
$ws.Cells.Item(27, 10).Value = 1.016202236447099
$ws.Cells.Item(27, 11).Value = 0.9961403935434965
$ws.Cells.Item(27, 12).Value = 0.9958433541799609
$ws.Cells.Item(27, 13).Value = 1.00602279531353
$ws.Cells.Item(27, 14).Value = 1.00602279531353
$ws.Cells.Item(27, 15).Value = 1.007008520396697
$ws.Cells.Item(27, 16).Value = 1.002728661390185
$ws.Cells.Item(27, 17).Value = 1.002728661390186
$ws.Cells.Item(27, 18).Value = 1.001081594428513
$ws.Cells.Item(27, 19).Value = 1.001081594428513
$ws.Cells.Item(27, 20).Value = 1.000325828036517
# Row 28
$ws.Cells.Item(28, 2).Value = "OffsetARD"
$ws.Cells.Item(28, 3).Value = 0.9972928526903513
$ws.Cells.Item(28, 4).Value = 1.002323749074494
$ws.Cells.Item(28, 5).Value = 0.99884873132929
$ws.Cells.Item(28, 6).Value = 0.9989387138255711
$ws.Cells.Item(28, 7).Value = 1.004278458247074
$ws.Cells.Item(28, 8).Value = 0.9989757454755455
$ws.Cells.Item(28, 9).Value = 0.99884873132929
$ws.Cells.Item(28, 10).Value = 1.004278458247074
$ws.Cells.Item(28, 11).Value = 0.99884873132929
$ws.Cells.Item(28, 12).Value = 0.9989757454755455
$ws.Cells.Item(28, 13).Value = 1.00162710186131
$ws.Cells.Item(28, 14).Value = 1.00162710186131
$ws.Cells.Item(28, 15).Value = 1.001859317599038
$ws.Cells.Item(28, 16).Value = 1.000700978350636
$ws.Cells.Item(28, 17).Value = 1.000700978350636
$ws.Cells.Item(28, 18).Value = 1.0002379165953
$ws.Cells.Item(28, 19).Value = 1.0002379165953
$ws.Cells.Item(28, 20).Value = 1.000109708440387
# Row 29
$ws.Cells.Item(29, 2).Value = "Gaussian Quadrature"
$ws.Cells.Item(29, 3).Value = 1.0049286431089
$ws.Cells.Item(29, 4).Value = 0.9959118157579536
$ws.Cells.Item(29, 5).Value = 1.001881450652574
$ws.Cells.Item(29, 6).Value = 1.001859489790283
$ws.Cells.Item(29, 7).Value = 0.9925539189243152
$ws.Cells.Item(29, 8).Value = 1.001850456768758
$ws.Cells.Item(29, 9).Value = 1.001881450652574
$ws.Cells.Item(29, 10).Value = 0.9925539189243152
$ws.Cells.Item(29, 11).Value = 1.001881450652574
$ws.Cells.Item(29, 12).Value = 1.001850456768758
$ws.Cells.Item(29, 13).Value = 0.9972021878465367
$ws.Cells.Item(29, 14).Value = 0.9972021878465367
$ws.Cells.Item(29, 15).Value = 0.996772063817009
$ws.Cells.Item(29, 16).Value = 0.9987619421152157
$ws.Cells.Item(29, 17).Value = 0.9987619421152156
$ws.Cells.Item(29, 18).Value = 0.9995418192495551
$ws.Cells.Item(29, 19).Value = 0.9995418192495551
$ws.Cells.Item(29, 20).Value = 0.9998309625004639
# Row 30
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Michael-CCHex"
$ws.Cells.Item(30, 3).Value = 0.9990092905082657
$ws.Cells.Item(30, 4).Value = 1.00070491312134
$ws.Cells.Item(30, 5).Value = 0.9997977969556445
$ws.Cells.Item(30, 6).Value = 0.9996858288555105
$ws.Cells.Item(30, 7).Value = 1.001215165832634
$ws.Cells.Item(30, 8).Value = 0.9996397701802237
$ws.Cells.Item(30, 9).Value = 0.9997977969556445
$ws.Cells.Item(30, 10).Value = 1.001215165832634
$ws.Cells.Item(30, 11).Value = 0.9997977969556445
$ws.Cells.Item(30, 12).Value = 0.9996397701802237
$ws.Cells.Item(30, 13).Value = 1.000427468006429
$ws.Cells.Item(30, 14).Value = 1.000427468006429
$ws.Cells.Item(30, 15).Value = 1.000519949711399
$ws.Cells.Item(30, 16).Value = 1.000217577656167
$ws.Cells.Item(30, 17).Value = 1.000217577656167
$ws.Cells.Item(30, 18).Value = 1.000112632481037
$ws.Cells.Item(30, 19).Value = 1.000112632481037
$ws.Cells.Item(30, 20).Value = 1.00000879424227
# Row 31
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Michael-SNHex"
$ws.Cells.Item(31, 3).Value = 0.9838213077783419
$ws.Cells.Item(31, 4).Value = 1.013181705608488
$ws.Cells.Item(31, 5).Value = 0.9941825897913593
$ws.Cells.Item(31, 6).Value = 0.9940175257118296
$ws.Cells.Item(31, 7).Value = 1.023868681282269
$ws.Cells.Item(31, 8).Value = 0.9939496282594635
$ws.Cells.Item(31, 9).Value = 0.9941825897913593
$ws.Cells.Item(31, 10).Value = 1.023868681282269
$ws.Cells.Item(31, 11).Value = 0.9941825897913593
$ws.Cells.Item(31, 12).Value = 0.9939496282594635
$ws.Cells.Item(31, 13).Value = 1.008909154770866
$ws.Cells.Item(31, 14).Value = 1.008909154770866
$ws.Cells.Item(31, 15).Value = 1.010333338383407
$ws.Cells.Item(31, 16).Value = 1.004000299777697
$ws.Cells.Item(31, 17).Value = 1.004000299777697
$ws.Cells.Item(31, 18).Value = 1.001545872281113
$ws.Cells.Item(31, 19).Value = 1.001545872281113
$ws.Cells.Item(31, 20).Value = 1.000503573071959

# New rows 30/31 need column A in the same bold/bordered style as the rest
# of the A column (copy the formatting down from the row above).
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
